# Create demo for data challenge.
# Append new test-condition rows (25-46) to "Test conditions.xlsx", Sheet1.
# Cell writes below are ordered to match the author's original entry order
# (new Test-id / Description strings were not typed strictly row-by-row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "20240524_094684"
$ws.Range("B25").Value = "Not moving"

$ws.Range("A26").Value = "20240524_094877"
$ws.Range("B26").Value = "Not moving"

$ws.Range("A27").Value = "20240524_100052"
$ws.Range("B27").Value = "Turning motor 1"

$ws.Range("B28").Value = "Turning motor 2"
$ws.Range("A28").Value = "20240524_101062"

$ws.Range("B29").Value = "Turning motor 3"

$ws.Range("B30").Value = "Turning motor 4"

$ws.Range("B31").Value = "Turning motor 5"

$ws.Range("A29").Value = "20240524_101487"

$ws.Range("A30").Value = "20240524_102066"

$ws.Range("A31").Value = "20240524_102736"

$ws.Range("A33").Value = "20240524_103973"

$ws.Range("A34").Value = "20240524_104453"

$ws.Range("A32").Value = "20240524_102301"
$ws.Range("B32").Value = "Turning motor 6"

$ws.Range("A35").Value = "20240524_104923"

$ws.Range("A36").Value = "20240524_105370"

$ws.Range("A37").Value = "20240524_105836"

$ws.Range("A38").Value = "20240524_110994"

$ws.Range("B33").Value = "Perform motor 1 fail"

$ws.Range("B34").Value = "Perform motor 2 fail"

$ws.Range("B35").Value = "Perform motor 3 fail"

$ws.Range("B36").Value = "Perform motor 4 fail"

$ws.Range("B37").Value = "Perform motor 5 fail"

$ws.Range("B38").Value = "Perform motor 321654 fail"

$ws.Range("A39").Value = "20240527_094865"
$ws.Range("B39").Value = "Transfer goods"
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0

$ws.Range("A40").Value = "20240527_100759"
$ws.Range("B40").Value = "Transfer goods"
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 1

$ws.Range("A41").Value = "20240527_101627"
$ws.Range("B41").Value = "Transfer goods"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0

$ws.Range("A42").Value = "20240527_102436"
$ws.Range("B42").Value = "Not moving"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 1

$ws.Range("A43").Value = "20240527_102919"
$ws.Range("B43").Value = "Not moving"
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1

$ws.Range("A44").Value = "20240527_103311"
$ws.Range("B44").Value = "Not moving"
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 1
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0

$ws.Range("A45").Value = "20240527_103690"
$ws.Range("B45").Value = "Moving one motor"
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 1

$ws.Range("A46").Value = "20240527_104247"
$ws.Range("B46").Value = "Moving one motor"
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

# --- View state: selection now sits on the newly-entered last row ------------
$null = $ws.Range("A46").Select()

# --- Page setup: portrait orientation (as captured in the saved file) --------
$ws.PageSetup.Orientation = 1
